$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date values in B1 and C1 (shift year from 2022 to 2023, same day/month)
$ws.Range("B1").Value = (Get-Date -Year 2023 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C1").Value = (Get-Date -Year 2023 -Month 4 -Day 14 -Hour 0 -Minute 0 -Second 0).Date

# Set explicit column widths (no longer best-fit / autofit)
$ws.Columns.Item(1).ColumnWidth = 12.8
$ws.Columns.Item(2).ColumnWidth = 14.67

# Update the selected cell shown when the sheet is opened
$ws.Range("B4").Select()
